$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.828.88"
$ws.Range("E2").Value = "  +2.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.815.29"
$ws.Range("E3").Value = "  +0.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "668.32"
$ws.Range("E5").Value = "  +6.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.11"
$ws.Range("E6").Value = "  +2.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.814.53"
$ws.Range("E7").Value = "  +0.84%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +1.34%  "

$ws.Range("E10").Value = "  +0.46%  "

$ws.Range("E11").Value = "  +2.39%  "

$ws.Range("E12").Value = "  +4.60%  "

$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.460.31"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.809.84"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.797.07"
$ws.Range("E17").Value = "  +2.04%  "

$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.65"
$ws.Range("E20").Value = "  +21.00%  "

$ws.Range("E21").Value = "  +0.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "475.32"
$ws.Range("E22").Value = "  +1.47%  "

$ws.Range("E23").Value = "  +1.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.24"

$ws.Range("E25").Value = "  -2.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.23"
$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.38"
$ws.Range("E27").Value = "  +3.50%  "

$ws.Range("E28").Value = "  -1.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.967.46"
$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.88"
$ws.Range("E31").Value = "  +7.91%  "

$ws.Range("E32").Value = "  +2.70%  "

$ws.Range("E33").Value = "  +1.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.75"
$ws.Range("E34").Value = "  +3.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.176"
$ws.Range("E35").Value = "  +5.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.18"
$ws.Range("E36").Value = "  +2.08%  "

$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.771.97"
$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("E41").Value = "  +3.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.968"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("E43").Value = "  +0.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.10"
$ws.Range("E44").Value = "  +9.36%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.63"
$ws.Range("E46").Value = "  +5.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.51"
$ws.Range("E47").Value = "  +3.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.01"
$ws.Range("E48").Value = "  +2.71%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.43"
$ws.Range("E49").Value = "  +4.39%  "

$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.300"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000291"
$ws.Range("E51").Value = "  +4.00%  "
